$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1826501.2
$ws.Range("I33").Value = 2130802.5
$ws.Range("J33").Value = 694
$ws.Range("K33").Value = 2130802.5
$ws.Range("L33").Value = 694
$ws.Range("M33").Value = -2130573.5
$ws.Range("N33").Value = -1152

$ws.Range("H40").Value = 23810710
$ws.Range("I40").Value = 1106.9375
$ws.Range("J40").Value = 38462776
$ws.Range("K40").Value = 1106.9375
$ws.Range("L40").Value = 38462776
$ws.Range("M40").Value = -931.9375
$ws.Range("N40").Value = -38463126

$ws.Range("H132").Value = 1467.76
$ws.Range("I132").Value = 879.0781
$ws.Range("J132").Value = 4892.8184
$ws.Range("K132").Value = 2637.2343
$ws.Range("L132").Value = 14678.4552
$ws.Range("M132").Value = -107.2343000000001
$ws.Range("N132").Value = -19738.4552

$ws.Range("H137").Value = 3126069.2
$ws.Range("I137").Value = 1191541.1
$ws.Range("J137").Value = 16667767
$ws.Range("K137").Value = 3574623.3
$ws.Range("L137").Value = 50003301
$ws.Range("M137").Value = -3572073.3
$ws.Range("N137").Value = -50008401

$ws.Range("H138").Value = 2396.25
$ws.Range("I138").Value = 2103.625
$ws.Range("J138").Value = 2786.4167
$ws.Range("K138").Value = 6310.875
$ws.Range("L138").Value = 8359.250100000001
$ws.Range("M138").Value = -1170.875
$ws.Range("N138").Value = -18639.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2964.4211
$ws.Range("I2").Value = 1671.5
$ws.Range("J2").Value = 9860
$ws.Range("K2").Value = 1671.5
$ws.Range("L2").Value = 9860
$ws.Range("M2").Value = -1558.5
$ws.Range("N2").Value = -10086

$ws.Range("H32").Value = 3543324.5
$ws.Range("I32").Value = 4731.949
$ws.Range("J32").Value = 20941404
$ws.Range("K32").Value = 4731.949
$ws.Range("L32").Value = 20941404
$ws.Range("M32").Value = -4444.949
$ws.Range("N32").Value = -20941978

$ws.Range("H116").Value = 2964.4211
$ws.Range("I116").Value = 1671.5
$ws.Range("J116").Value = 9860
$ws.Range("K116").Value = 1671.5
$ws.Range("L116").Value = 9860
$ws.Range("M116").Value = 622.5
$ws.Range("N116").Value = -14448

$ws.Range("H132").Value = 97425.445
$ws.Range("I132").Value = 103165.81
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 309497.43
$ws.Range("L132").Value = 10998
$ws.Range("M132").Value = -306967.43
$ws.Range("N132").Value = -16058

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2964.4211
$ws.Range("I3").Value = 1671.5
$ws.Range("J3").Value = 9860
$ws.Range("K3").Value = 1671.5
$ws.Range("L3").Value = 9860
$ws.Range("M3").Value = -1557.5
$ws.Range("N3").Value = -10088

$ws.Range("H99").Value = 1241.1111
$ws.Range("I99").Value = 896
$ws.Range("J99").Value = 2966.6667
$ws.Range("K99").Value = 896
$ws.Range("L99").Value = 2966.6667
$ws.Range("M99").Value = 602
$ws.Range("N99").Value = -5962.6667

$ws.Range("H107").Value = 10726.375
$ws.Range("I107").Value = 10830.143
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 10830.143
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -8910.143
$ws.Range("N107").Value = -13840

$ws.Range("H134").Value = 121254.92
$ws.Range("I134").Value = 131577.39
$ws.Range("J134").Value = 2546.5
$ws.Range("K134").Value = 394732.17
$ws.Range("L134").Value = 7639.5
$ws.Range("M134").Value = -392197.17
$ws.Range("N134").Value = -12709.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1846.7646
$ws.Range("I31").Value = 1563.7894
$ws.Range("J31").Value = 2205.2
$ws.Range("K31").Value = 1563.7894
$ws.Range("L31").Value = 2205.2
$ws.Range("M31").Value = -1268.7894
$ws.Range("N31").Value = -2795.2

$ws.Range("H34").Value = 1846.7646
$ws.Range("I34").Value = 1563.7894
$ws.Range("J34").Value = 2205.2
$ws.Range("K34").Value = 1563.7894
$ws.Range("L34").Value = 2205.2
$ws.Range("M34").Value = -1361.7894
$ws.Range("N34").Value = -2609.2

$ws.Range("H132").Value = 3114.0344
$ws.Range("I132").Value = 2555
$ws.Range("J132").Value = 4871
$ws.Range("K132").Value = 7665
$ws.Range("L132").Value = 14613
$ws.Range("M132").Value = -5135
$ws.Range("N132").Value = -19673

$ws.Range("H134").Value = 4209.281
$ws.Range("I134").Value = 4641.9585
$ws.Range("J134").Value = 1901.6666
$ws.Range("K134").Value = 13925.8755
$ws.Range("L134").Value = 5704.9998
$ws.Range("M134").Value = -11390.8755
$ws.Range("N134").Value = -10774.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 982.098
$ws.Range("J48").Value = 1022.5833
$ws.Range("L48").Value = 3067.7499
$ws.Range("N48").Value = -3567.7499

$ws.Range("H122").Value = 48135.945
$ws.Range("I122").Value = 345.45456
$ws.Range("K122").Value = 3109.09104
$ws.Range("M122").Value = -659.0910400000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2607
$ws.Range("I132").Value = 2084.182
$ws.Range("J132").Value = 2926.5
$ws.Range("K132").Value = 6252.545999999999
$ws.Range("L132").Value = 8779.5
$ws.Range("M132").Value = -3722.545999999999
$ws.Range("N132").Value = -13839.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1771.4286
$ws.Range("I7").Value = 1700
$ws.Range("K7").Value = 1700
$ws.Range("M7").Value = -1588

$ws.Range("H126").Value = 1771.4286
$ws.Range("I126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("M126").Value = -2630

$ws.Range("H136").Value = 2158.0908
$ws.Range("I136").Value = 1840.68
$ws.Range("J136").Value = 3150
$ws.Range("K136").Value = 5522.04
$ws.Range("L136").Value = 9450
$ws.Range("M136").Value = -2972.04
$ws.Range("N136").Value = -14550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1856.5
$ws.Range("I96").Value = 1722.5
$ws.Range("K96").Value = 1722.5
$ws.Range("M96").Value = -349.5

$ws.Range("H122").Value = 3792.4285
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4232.6665
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 12697.9995
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -17597.9995

$ws.Range("H132").Value = 2227.0625
$ws.Range("I132").Value = 2075.3462
$ws.Range("J132").Value = 2884.5
$ws.Range("K132").Value = 6226.0386
$ws.Range("L132").Value = 8653.5
$ws.Range("M132").Value = -3696.0386
$ws.Range("N132").Value = -13713.5
